$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.322.75"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'2.822.76"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'361.16"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'111.29"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").Value = "'0.566"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").Value = "'40.58"
$ws.Range("E10").Value = "  -5.58%  "
$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'19.82"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "'3.262.43"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "'2.872.48"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "'0.921"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "'52.092.88"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'7.48"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "'3.14"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "'13.34"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "'0.0₃0995"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "'273.46"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'70.12"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "'2.81"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'26.81"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'10.26"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "'0.0477"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").Value = "'52.23"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "'34.44"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").Value = "'5.83"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "'5.49"
$ws.Range("E35").Value = "  +10.62%  "
$ws.Range("D36").Value = "'0.0849"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "'2.03"
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("D40").Value = "'18.14"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "'2.55"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("D43").Value = "'125.20"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "'22.41"
$ws.Range("E45").Value = "  -6.19%  "
$ws.Range("D46").Value = "'2.071.11"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'3.30"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").Value = "'2.35"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("D49").Value = "'5.85"
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("D50").Value = "'0.951"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "'9.11"
$ws.Range("E51").Value = "  +0.73%  "
